# Update the NATMI LR-pair (Wnt11-Fzd8) TPM-derived metrics on Sheet1.
# Columns A-F (cluster/ligand/receptor labels, cell counts, detection
# rates) are unchanged; columns G-T (expression/specificity/edge-weight
# values derived from the refreshed TPM input) are updated in place for
# every data row (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,7).Value = 0.270956
$ws.Cells.Item(2,8).Value = 0.812868
$ws.Cells.Item(2,9).Value = 0.05752762181187367
$ws.Cells.Item(2,10).Value = 0.05752762181187366
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 3.024701
$ws.Cells.Item(2,14).Value = 9.074103000000001
$ws.Cells.Item(2,15).Value = 0.1596375877334842
$ws.Cells.Item(2,16).Value = 0.1596375877334843
$ws.Cells.Item(2,17).Value = 0.8195608841560001
$ws.Cells.Item(2,18).Value = 7.376047957404001
$ws.Cells.Item(2,19).Value = 0.009183570774091683
$ws.Cells.Item(2,20).Value = 0.009183570774091685
$ws.Cells.Item(3,7).Value = 0.270956
$ws.Cells.Item(3,8).Value = 0.812868
$ws.Cells.Item(3,9).Value = 0.05752762181187367
$ws.Cells.Item(3,10).Value = 0.05752762181187366
$ws.Cells.Item(3,15).Value = 0.6072559333217162
$ws.Cells.Item(3,16).Value = 0.6072559333217163
$ws.Cells.Item(3,17).Value = 3.117581621522667
$ws.Cells.Item(3,18).Value = 28.058234593704
$ws.Cells.Item(3,19).Value = 0.03493398967514806
$ws.Cells.Item(3,20).Value = 0.03493398967514807
$ws.Cells.Item(4,7).Value = 0.270956
$ws.Cells.Item(4,8).Value = 0.812868
$ws.Cells.Item(4,9).Value = 0.05752762181187367
$ws.Cells.Item(4,10).Value = 0.05752762181187366
$ws.Cells.Item(4,13).Value = 4.368554666666666
$ws.Cells.Item(4,14).Value = 13.105664
$ws.Cells.Item(4,15).Value = 0.2305634602787257
$ws.Cells.Item(4,16).Value = 0.2305634602787257
$ws.Cells.Item(4,17).Value = 1.183686098261333
$ws.Cells.Item(4,18).Value = 10.653174884352
$ws.Cells.Item(4,19).Value = 0.01326376754655149
$ws.Cells.Item(4,20).Value = 0.01326376754655149
$ws.Cells.Item(5,7).Value = 0.270956
$ws.Cells.Item(5,8).Value = 0.812868
$ws.Cells.Item(5,9).Value = 0.05752762181187367
$ws.Cells.Item(5,10).Value = 0.05752762181187366
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.04818333333333333
$ws.Cells.Item(5,14).Value = 0.14455
$ws.Cells.Item(5,15).Value = 0.002543018666073676
$ws.Cells.Item(5,16).Value = 0.002543018666073677
$ws.Cells.Item(5,17).Value = 0.01305556326666667
$ws.Cells.Item(5,18).Value = 0.1175000694
$ws.Cells.Item(5,19).Value = 0.0001462938160824219
$ws.Cells.Item(5,20).Value = 0.0001462938160824219
$ws.Cells.Item(6,9).Value = 0.7443710555244437
$ws.Cells.Item(6,10).Value = 0.7443710555244437
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.024701
$ws.Cells.Item(6,14).Value = 9.074103000000001
$ws.Cells.Item(6,15).Value = 0.1596375877334842
$ws.Cells.Item(6,16).Value = 0.1596375877334843
$ws.Cells.Item(6,17).Value = 10.60459968953267
$ws.Cells.Item(6,18).Value = 95.44139720579402
$ws.Cells.Item(6,19).Value = 0.1188295996825496
$ws.Cells.Item(6,20).Value = 0.1188295996825497
$ws.Cells.Item(7,9).Value = 0.7443710555244437
$ws.Cells.Item(7,10).Value = 0.7443710555244437
$ws.Cells.Item(7,15).Value = 0.6072559333217162
$ws.Cells.Item(7,16).Value = 0.6072559333217163
$ws.Cells.Item(7,19).Value = 0.4520237400601671
$ws.Cells.Item(7,20).Value = 0.4520237400601672
$ws.Cells.Item(8,9).Value = 0.7443710555244437
$ws.Cells.Item(8,10).Value = 0.7443710555244437
$ws.Cells.Item(8,13).Value = 4.368554666666666
$ws.Cells.Item(8,14).Value = 13.105664
$ws.Cells.Item(8,15).Value = 0.2305634602787257
$ws.Cells.Item(8,16).Value = 0.2305634602787257
$ws.Cells.Item(8,17).Value = 15.31614974896356
$ws.Cells.Item(8,18).Value = 137.845347740672
$ws.Cells.Item(8,19).Value = 0.1716247662930432
$ws.Cells.Item(8,20).Value = 0.1716247662930432
$ws.Cells.Item(9,9).Value = 0.7443710555244437
$ws.Cells.Item(9,10).Value = 0.7443710555244437
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.04818333333333333
$ws.Cells.Item(9,14).Value = 0.14455
$ws.Cells.Item(9,15).Value = 0.002543018666073676
$ws.Cells.Item(9,16).Value = 0.002543018666073677
$ws.Cells.Item(9,17).Value = 0.1689307345444445
$ws.Cells.Item(9,18).Value = 1.5203766109
$ws.Cells.Item(9,19).Value = 0.001892949488683625
$ws.Cells.Item(9,20).Value = 0.001892949488683626
$ws.Cells.Item(10,7).Value = 0.3861913333333333
$ws.Cells.Item(10,8).Value = 1.158574
$ws.Cells.Item(10,9).Value = 0.08199364092702592
$ws.Cells.Item(10,10).Value = 0.08199364092702592
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.024701
$ws.Cells.Item(10,14).Value = 9.074103000000001
$ws.Cells.Item(10,15).Value = 0.1596375877334842
$ws.Cells.Item(10,16).Value = 0.1596375877334843
$ws.Cells.Item(10,17).Value = 1.168113312124667
$ws.Cells.Item(10,18).Value = 10.513019809122
$ws.Cells.Item(10,19).Value = 0.01308926704707591
$ws.Cells.Item(10,20).Value = 0.01308926704707591
$ws.Cells.Item(11,7).Value = 0.3861913333333333
$ws.Cells.Item(11,8).Value = 1.158574
$ws.Cells.Item(11,9).Value = 0.08199364092702592
$ws.Cells.Item(11,10).Value = 0.08199364092702592
$ws.Cells.Item(11,15).Value = 0.6072559333217162
$ws.Cells.Item(11,16).Value = 0.6072559333217163
$ws.Cells.Item(11,17).Value = 4.443463157085778
$ws.Cells.Item(11,18).Value = 39.991168413772
$ws.Cells.Item(11,19).Value = 0.04979112494758679
$ws.Cells.Item(11,20).Value = 0.04979112494758681
$ws.Cells.Item(12,7).Value = 0.3861913333333333
$ws.Cells.Item(12,8).Value = 1.158574
$ws.Cells.Item(12,9).Value = 0.08199364092702592
$ws.Cells.Item(12,10).Value = 0.08199364092702592
$ws.Cells.Item(12,13).Value = 4.368554666666666
$ws.Cells.Item(12,14).Value = 13.105664
$ws.Cells.Item(12,15).Value = 0.2305634602787257
$ws.Cells.Item(12,16).Value = 0.2305634602787257
$ws.Cells.Item(12,17).Value = 1.687097951459555
$ws.Cells.Item(12,18).Value = 15.183881563136
$ws.Cells.Item(12,19).Value = 0.01890473757298644
$ws.Cells.Item(12,20).Value = 0.01890473757298644
$ws.Cells.Item(13,7).Value = 0.3861913333333333
$ws.Cells.Item(13,8).Value = 1.158574
$ws.Cells.Item(13,9).Value = 0.08199364092702592
$ws.Cells.Item(13,10).Value = 0.08199364092702592
$ws.Cells.Item(13,11).Value = 2
$ws.Cells.Item(13,12).Value = 0.6666666666666666
$ws.Cells.Item(13,13).Value = 0.04818333333333333
$ws.Cells.Item(13,14).Value = 0.14455
$ws.Cells.Item(13,15).Value = 0.002543018666073676
$ws.Cells.Item(13,16).Value = 0.002543018666073677
$ws.Cells.Item(13,17).Value = 0.01860798574444444
$ws.Cells.Item(13,18).Value = 0.1674718717
$ws.Cells.Item(13,19).Value = 0.0002085113593767695
$ws.Cells.Item(13,20).Value = 0.0002085113593767695
$ws.Cells.Item(14,7).Value = 0.546869
$ws.Cells.Item(14,8).Value = 1.640607
$ws.Cells.Item(14,9).Value = 0.1161076817366566
$ws.Cells.Item(14,10).Value = 0.1161076817366566
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 3.024701
$ws.Cells.Item(14,14).Value = 9.074103000000001
$ws.Cells.Item(14,15).Value = 0.1596375877334842
$ws.Cells.Item(14,16).Value = 0.1596375877334843
$ws.Cells.Item(14,17).Value = 1.654115211169
$ws.Cells.Item(14,18).Value = 14.887036900521
$ws.Cells.Item(14,19).Value = 0.01853515022976699
$ws.Cells.Item(14,20).Value = 0.01853515022976699
$ws.Cells.Item(15,7).Value = 0.546869
$ws.Cells.Item(15,8).Value = 1.640607
$ws.Cells.Item(15,9).Value = 0.1161076817366566
$ws.Cells.Item(15,10).Value = 0.1161076817366566
$ws.Cells.Item(15,15).Value = 0.6072559333217162
$ws.Cells.Item(15,16).Value = 0.6072559333217163
$ws.Cells.Item(15,17).Value = 6.292197787760667
$ws.Cells.Item(15,18).Value = 56.62978008984601
$ws.Cells.Item(15,19).Value = 0.07050707863881421
$ws.Cells.Item(15,20).Value = 0.07050707863881422
$ws.Cells.Item(16,7).Value = 0.546869
$ws.Cells.Item(16,8).Value = 1.640607
$ws.Cells.Item(16,9).Value = 0.1161076817366566
$ws.Cells.Item(16,10).Value = 0.1161076817366566
$ws.Cells.Item(16,13).Value = 4.368554666666666
$ws.Cells.Item(16,14).Value = 13.105664
$ws.Cells.Item(16,15).Value = 0.2305634602787257
$ws.Cells.Item(16,16).Value = 0.2305634602787257
$ws.Cells.Item(16,17).Value = 2.389027122005333
$ws.Cells.Item(16,18).Value = 21.501244098048
$ws.Cells.Item(16,19).Value = 0.02677018886614456
$ws.Cells.Item(16,20).Value = 0.02677018886614456
$ws.Cells.Item(17,7).Value = 0.546869
$ws.Cells.Item(17,8).Value = 1.640607
$ws.Cells.Item(17,9).Value = 0.1161076817366566
$ws.Cells.Item(17,10).Value = 0.1161076817366566
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.04818333333333333
$ws.Cells.Item(17,14).Value = 0.14455
$ws.Cells.Item(17,15).Value = 0.002543018666073676
$ws.Cells.Item(17,16).Value = 0.002543018666073677
$ws.Cells.Item(17,17).Value = 0.02634997131666667
$ws.Cells.Item(17,18).Value = 0.23714974185
$ws.Cells.Item(17,19).Value = 0.0002952640019308595
$ws.Cells.Item(17,20).Value = 0.0002952640019308596
